$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything that changes (keep A1 exactly as-is: the blank corner cell)
# so the shared-string table rebuilds fresh, in final row-major first-use order.
$ws.Range("B1:O1").Clear()
$ws.Range("A2:O11").Clear()

# Header row (B1:K1) -- Europe columns dropped, remaining groups shift left
$ws.Range("B1").Value = "`$ bold('All')"
$ws.Range("C1").Value = "Millionaires"
$ws.Range("D1").Value = "Japan Non-voters"
$ws.Range("E1").Value = "Japan Left"
$ws.Range("F1").Value = "Japan Center/Right"
$ws.Range("G1").Value = "Saudi Arabia"
$ws.Range("H1").Value = "Saudi citizens"
$ws.Range("I1").Value = "U.S. Non-voters"
$ws.Range("J1").Value = "U.S. Harris"
$ws.Range("K1").Value = "U.S. Trump"

# Data rows 2-11: policy name + per-group shares, final recomputed data,
# rows resorted by the updated "All" share (column B), descending
$ws.Range("A2").Value = "Minimum tax of 2% on billionaires'`nwealth, in voluntary countries"
$ws.Range("B2").Value = 0.809187541057965
$ws.Range("C2").Value = 0.643085962793994
$ws.Range("D2").Value = 0.819942754959698
$ws.Range("E2").Value = 0.873056146222153
$ws.Range("F2").Value = 0.773791371449113
$ws.Range("G2").Value = 0.859553623962162
$ws.Range("H2").Value = 0.852102558492413
$ws.Range("I2").Value = 0.808151291496163
$ws.Range("J2").Value = 0.931123356175572
$ws.Range("K2").Value = 0.563794665192873

$ws.Range("A3").Value = "Bridgetown initiative: MDBs expanding sustainable`ninvestments in LICs, and at lower interest rates"
$ws.Range("B3").Value = 0.793366965036613
$ws.Range("C3").Value = 0.735290211926387
$ws.Range("D3").Value = 0.767899521650745
$ws.Range("E3").Value = 0.856662637816024
$ws.Range("F3").Value = 0.80190982641243
$ws.Range("G3").Value = 0.870851551659494
$ws.Range("H3").Value = 0.86802251791901
$ws.Range("I3").Value = 0.779977469932862
$ws.Range("J3").Value = 0.91765670226024
$ws.Range("K3").Value = 0.51808890832465

$ws.Range("A4").Value = "L&D: Developed countries financing a fund to help`nvulnerable countries cope with climate Loss and damage"
$ws.Range("B4").Value = 0.748295871658231
$ws.Range("C4").Value = 0.643308316798942
$ws.Range("D4").Value = 0.719352832772282
$ws.Range("E4").Value = 0.773521954970104
$ws.Range("F4").Value = 0.724550091724308
$ws.Range("G4").Value = 0.894499688071821
$ws.Range("H4").Value = 0.900586001863265
$ws.Range("I4").Value = 0.746047610561908
$ws.Range("J4").Value = 0.884011688881496
$ws.Range("K4").Value = 0.453402045694884

$ws.Range("A5").Value = "International levy on shipping carbon emissions,`nreturned to countries based on population"
$ws.Range("B5").Value = 0.699644803977167
$ws.Range("C5").Value = 0.596447568379874
$ws.Range("D5").Value = 0.565647078769302
$ws.Range("E5").Value = 0.635964936582196
$ws.Range("F5").Value = 0.585971465473506
$ws.Range("G5").Value = 0.814900578705803
$ws.Range("H5").Value = 0.842045082491013
$ws.Range("I5").Value = 0.684753896310072
$ws.Range("J5").Value = 0.86203143431525
$ws.Range("K5").Value = 0.471193151072884

$ws.Range("A6").Value = "At least 0.7% of developed countries' GDP in foreign aid"
$ws.Range("B6").Value = 0.698715666285492
$ws.Range("C6").Value = 0.62563884053293
$ws.Range("D6").Value = 0.564360282838369
$ws.Range("E6").Value = 0.656346841758792
$ws.Range("F6").Value = 0.657746953380401
$ws.Range("G6").Value = 0.863575793802146
$ws.Range("H6").Value = 0.871078497471275
$ws.Range("I6").Value = 0.681660037168233
$ws.Range("J6").Value = 0.839986172033593
$ws.Range("K6").Value = 0.458847635276405

$ws.Range("A7").Value = "Debt relief for vulnerable countries, suspending`npayments until they are more able to repay"
$ws.Range("B7").Value = 0.696886644817168
$ws.Range("C7").Value = 0.546116992829272
$ws.Range("D7").Value = 0.686942886007181
$ws.Range("E7").Value = 0.742202702288568
$ws.Range("F7").Value = 0.658450500259319
$ws.Range("G7").Value = 0.883360143018035
$ws.Range("H7").Value = 0.898514990349396
$ws.Range("I7").Value = 0.74293484896282
$ws.Range("J7").Value = 0.806053340993667
$ws.Range("K7").Value = 0.468557448120917

$ws.Range("A8").Value = "Expand Security Council to new permanent members (e.g.`nIndia, Brazil, African Union), restrict veto use"
$ws.Range("B8").Value = 0.694962460493881
$ws.Range("C8").Value = 0.618548578698709
$ws.Range("D8").Value = 0.643460694813484
$ws.Range("E8").Value = 0.715145780417729
$ws.Range("F8").Value = 0.682980405531271
$ws.Range("G8").Value = 0.836791104476278
$ws.Range("H8").Value = 0.854225312098256
$ws.Range("I8").Value = 0.6953146471654
$ws.Range("J8").Value = 0.861220065079599
$ws.Range("K8").Value = 0.454492739666233

$ws.Range("A9").Value = "NCQG: Developing countries providing `$300 bn a`nyear in climate finance for developing countries"
$ws.Range("B9").Value = 0.683174203642518
$ws.Range("C9").Value = 0.553445646083967
$ws.Range("D9").Value = 0.514556584884823
$ws.Range("E9").Value = 0.656047851343109
$ws.Range("F9").Value = 0.60344269229444
$ws.Range("G9").Value = 0.858825470219245
$ws.Range("H9").Value = 0.868267260484353
$ws.Range("I9").Value = 0.669358374420139
$ws.Range("J9").Value = 0.830923870403501
$ws.Range("K9").Value = 0.343601913157298

$ws.Range("A10").Value = "Raise global minimum tax on profit from 15% to 35%,`nallocating revenues to countries based on sales"
$ws.Range("B10").Value = 0.682119272063001
$ws.Range("C10").Value = 0.606189372320828
$ws.Range("D10").Value = 0.698572259799889
$ws.Range("E10").Value = 0.770942809545018
$ws.Range("F10").Value = 0.710087508866346
$ws.Range("G10").Value = 0.770097780765828
$ws.Range("H10").Value = 0.763530693008704
$ws.Range("I10").Value = 0.661649896053245
$ws.Range("J10").Value = 0.841576010118254
$ws.Range("K10").Value = 0.450190602157014

$ws.Range("A11").Value = "International levy on aviation carbon emissions, raising`nprices by 30%, returned to countries based on population"
$ws.Range("B11").Value = 0.526319576644156
$ws.Range("C11").Value = 0.458534242678225
$ws.Range("D11").Value = 0.432517119215557
$ws.Range("E11").Value = 0.485374406010005
$ws.Range("F11").Value = 0.474756209508454
$ws.Range("G11").Value = 0.695189510992556
$ws.Range("H11").Value = 0.729294481740513
$ws.Range("I11").Value = 0.475175961400892
$ws.Range("J11").Value = 0.673509672793868
$ws.Range("K11").Value = 0.338595491710254
